{"js": "// Figure \"Monarch\" caption: append description of the one-break-model\n// alternate fit, and relocate the (hidden, cursor-tracking) \"_GoBack\"\n// bookmark so it now collapses right before the newly-added clause\n// (\"phases B and C, is given by the black dashed line\") instead of\n// sitting at the start of the paragraph.\n\nconst body = context.document.body;\n\n// 1) Locate the end of the existing sentence we are extending.\nconst anchorResults = body.search(\n  \"in the transition from phase A to phase B\",\n  { matchCase: true, matchWholeWord: false }\n);\nanchorResults.load(\"text\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error('Could not find the anchor text to extend the Monarch caption.');\n}\n\nconst anchor = anchorResults.items[0];\n\n// 2) Append the new sentence describing the alternate (one break model) fit.\nanchor.insertText(\n  \". An alternate fit associated with a one break model that combine phases B and C, is given by the black dashed line.\",\n  \"After\"\n);\nawait context.sync();\n\n// 3) Remove the \"_GoBack\" bookmark from its old (collapsed) location near\n//    the start of the paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) Re-insert \"_GoBack\" collapsed immediately before \"phases B and C...\",\n//    matching its new position in the edited caption.\nconst newSpotResults = body.search(\n  \"phases B and C, is given by the black dashed line\",\n  { matchCase: true, matchWholeWord: false }\n);\nnewSpotResults.load(\"text\");\nawait context.sync();\n\nif (newSpotResults.items.length === 0) {\n  throw new Error('Could not find the new bookmark anchor text.');\n}\n\nconst newSpotStart = newSpotResults.items[0].getRange(\"Start\");\nnewSpotStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Figure \"Monarch\" caption: append description of the one-break-model\n# alternate fit, and relocate the (hidden, cursor-tracking) \"_GoBack\"\n# bookmark so it now collapses right before the newly-added clause\n# (\"phases B and C, is given by the black dashed line\") instead of\n# sitting at the start of the paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Find the end of the existing sentence we are extending.\n$anchor = $d.Content\n$found = $anchor.Find.Execute(\"in the transition from phase A to phase B\")\nif (-not $found) {\n    throw \"Could not find the anchor text to extend the Monarch caption.\"\n}\n\n# 2) Collapse to the end of that text and append the new sentence describing\n#    the alternate (one break model) fit.\n$anchor.Collapse(0)  # wdCollapseEnd\n$anchor.InsertAfter(\". An alternate fit associated with a one break model that combine phases B and C, is given by the black dashed line.\")\n\n# 3) Remove the \"_GoBack\" bookmark from its old (collapsed) location near\n#    the start of the paragraph.\n$goBack = $d.Bookmarks.Item(\"_GoBack\")\n$goBack.Delete()\n\n# 4) Re-insert \"_GoBack\" collapsed immediately before \"phases B and C...\",\n#    matching its new position in the edited caption.\n$newSpot = $d.Content\n$found2 = $newSpot.Find.Execute(\"phases B and C, is given by the black dashed line\")\nif (-not $found2) {\n    throw \"Could not find the new bookmark anchor text.\"\n}\n$newSpot.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $newSpot)\n"}
